# Auto update Excel log
# Append 7 new PRESENCE_DETECTED log rows (rows 61-67) to the "mmWave" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$rows = @(
    @{ Row = 61; Date = "2026-02-01"; Time = "17:45:21"; Hour = "17:00"; Location = "Living Room"; Value = "PRESENCE_DETECTED"; Status = "Active" },
    @{ Row = 62; Date = "2026-02-01"; Time = "17:45:31"; Hour = "17:00"; Location = "Living Room"; Value = "PRESENCE_DETECTED"; Status = "Active" },
    @{ Row = 63; Date = "2026-02-01"; Time = "17:47:58"; Hour = "17:00"; Location = "Living Room"; Value = "PRESENCE_DETECTED"; Status = "Active" },
    @{ Row = 64; Date = "2026-02-01"; Time = "17:48:08"; Hour = "17:00"; Location = "Living Room"; Value = "PRESENCE_DETECTED"; Status = "Active" },
    @{ Row = 65; Date = "2026-02-01"; Time = "17:48:19"; Hour = "17:00"; Location = "Living Room"; Value = "PRESENCE_DETECTED"; Status = "Active" },
    @{ Row = 66; Date = "2026-02-01"; Time = "17:48:29"; Hour = "17:00"; Location = "Living Room"; Value = "PRESENCE_DETECTED"; Status = "Active" },
    @{ Row = 67; Date = "2026-02-01"; Time = "17:48:40"; Hour = "17:00"; Location = "Living Room"; Value = "PRESENCE_DETECTED"; Status = "Active" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column A holds a date-like string ("2026-02-01"). Format the cell as
    # Text first so Excel keeps it as a literal string instead of silently
    # converting it to a date serial number.
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.Date

    $ws.Cells.Item($rowNum, 2).Value = $r.Time
    $ws.Cells.Item($rowNum, 3).Value = $r.Hour
    $ws.Cells.Item($rowNum, 4).Value = $r.Location
    $ws.Cells.Item($rowNum, 5).Value = $r.Value
    $ws.Cells.Item($rowNum, 6).Value = $r.Status
}
